$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 15, shifting existing rows 15-68 down to 19-72.
# This automatically preserves all formatting/content of the old rows in their new positions.
$ws.Range("A15:T18").EntireRow.Insert()

# New row 15: Helena, Especial
$ws.Cells.Item(15, 1).Value = 9
$ws.Cells.Item(15, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(15, 3).Value = "Metropolitana"
$ws.Cells.Item(15, 4).Value = 44558
$ws.Cells.Item(15, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value = 13
$ws.Cells.Item(15, 6).Value = "Fruta"
$ws.Cells.Item(15, 7).Value = 100103
$ws.Cells.Item(15, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(15, 9).Value = 100103003
$ws.Cells.Item(15, 10).Value = "Damasco"
$ws.Cells.Item(15, 11).Value = "Helena"
$ws.Cells.Item(15, 12).Value = "Especial"
$ws.Cells.Item(15, 13).Value = 330
$ws.Cells.Item(15, 14).Value = 16000
$ws.Cells.Item(15, 15).Value = 16000
$ws.Cells.Item(15, 16).Value = 16000
$ws.Cells.Item(15, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(15, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(15, 19).Value = 1000
$ws.Cells.Item(15, 20).Value = 16

# New row 16: Helena, Primera
$ws.Cells.Item(16, 1).Value = 9
$ws.Cells.Item(16, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44558
$ws.Cells.Item(16, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100103
$ws.Cells.Item(16, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(16, 9).Value = 100103003
$ws.Cells.Item(16, 10).Value = "Damasco"
$ws.Cells.Item(16, 11).Value = "Helena"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 350
$ws.Cells.Item(16, 14).Value = 12800
$ws.Cells.Item(16, 15).Value = 12800
$ws.Cells.Item(16, 16).Value = 12800
$ws.Cells.Item(16, 17).Value = "`$/caja 16 kilos"
$ws.Cells.Item(16, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(16, 19).Value = 800
$ws.Cells.Item(16, 20).Value = 16

# New row 17: Patterson, Primera
$ws.Cells.Item(17, 1).Value = 9
$ws.Cells.Item(17, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(17, 3).Value = "Metropolitana"
$ws.Cells.Item(17, 4).Value = 44558
$ws.Cells.Item(17, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(17, 5).Value = 13
$ws.Cells.Item(17, 6).Value = "Fruta"
$ws.Cells.Item(17, 7).Value = 100103
$ws.Cells.Item(17, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(17, 9).Value = 100103003
$ws.Cells.Item(17, 10).Value = "Damasco"
$ws.Cells.Item(17, 11).Value = "Patterson"
$ws.Cells.Item(17, 12).Value = "Primera"
$ws.Cells.Item(17, 13).Value = 410
$ws.Cells.Item(17, 14).Value = 14000
$ws.Cells.Item(17, 15).Value = 14000
$ws.Cells.Item(17, 16).Value = 14000
$ws.Cells.Item(17, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(17, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(17, 19).Value = 778
$ws.Cells.Item(17, 20).Value = 18

# New row 18: Patterson, Segunda
$ws.Cells.Item(18, 1).Value = 9
$ws.Cells.Item(18, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 44558
$ws.Cells.Item(18, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100103
$ws.Cells.Item(18, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(18, 9).Value = 100103003
$ws.Cells.Item(18, 10).Value = "Damasco"
$ws.Cells.Item(18, 11).Value = "Patterson"
$ws.Cells.Item(18, 12).Value = "Segunda"
$ws.Cells.Item(18, 13).Value = 380
$ws.Cells.Item(18, 14).Value = 12000
$ws.Cells.Item(18, 15).Value = 12000
$ws.Cells.Item(18, 16).Value = 12000
$ws.Cells.Item(18, 17).Value = "`$/caja 18 kilos"
$ws.Cells.Item(18, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(18, 19).Value = 667
$ws.Cells.Item(18, 20).Value = 18

$ws.Range("A1").Select()
